$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 213, shifting rows 213:281 down to 214:282
$ws.Rows.Item(213).Insert()

# Populate the new row 213 with its data (categorical columns copied
# from the row that was shifted down, numeric columns updated)
$ws.Cells.Item(213, 1).Value = 9
$ws.Cells.Item(213, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(213, 3).Value = "Metropolitana"
$ws.Cells.Item(213, 4).Value = 44841
$ws.Cells.Item(213, 5).Value = 13
$ws.Cells.Item(213, 6).Value = 100112026
$ws.Cells.Item(213, 7).Value = "Haba"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Primera"
$ws.Cells.Item(213, 10).Value = 114
$ws.Cells.Item(213, 11).Value = 12000
$ws.Cells.Item(213, 12).Value = 12000
$ws.Cells.Item(213, 13).Value = 12000
$ws.Cells.Item(213, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(213, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(213, 16).Value = 480
$ws.Cells.Item(213, 17).Value = 25
$ws.Cells.Item(213, 18).Value = "Hortaliza"
